# Add two new patient rows (book_appointments) to the PATIENT sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - kushal@practo.com
$ws.Range("A2").Value = "kushal@practo.com"
$ws.Range("B2").Value = "passwd123"
$ws.Range("C2").Value = "MALE"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 9087654321

# Row 3 - abhishek@practo.com
$ws.Range("A3").Value = "abhishek@practo.com"
$ws.Range("B3").Value = "passwd123"
$ws.Range("C3").Value = "MALE"
$ws.Range("D3").Value = 20
# Phone number stored as text on this row (leading apostrophe forces text, like typing it in Excel)
$ws.Range("E3").Value = "'9087654321"
